$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4 for the new requirement, pushing everything
# else (old rows 4-9) down by one (new rows 5-10).
$ws.Rows("4").Insert()

# Row 4: new requirement - "Create API for adding datasources, tags and
# information. Used for admin web app" with wrap style (matches rows that
# already wrap, e.g. the "Create web application..." row) and B=5.
$ws.Range("A4").Value = "Create API for adding datasources, tags and information. Used for admin web app"
$ws.Range("A4").WrapText = $true
$ws.Rows("4").RowHeight = 28.8
$ws.Range("B4").Value = 5

# Row 5: "Create web application..." (unchanged text/style), days now 15.
$ws.Range("B5").Value = 15

# Row 6: "Web Services..." (unchanged text/style), days now 15.
$ws.Range("B6").Value = 15

# Row 7: "Create Android application" (unchanged text/style), days now 25.
$ws.Range("B7").Value = 25

# Rows 8 & 9 swap places: "add offline capabilities..." moves up to row 8
# (taking on the wrap style/row height previously used by "provide Sabisu
# authentication..."), and "provide Sabisu authentication..." moves down to
# row 9 (taking on the plain style previously used by "add offline...").
$ws.Range("A8").Value = "add offline capabilities for storing tag locations and most recent data"
$ws.Range("A8").WrapText = $true
$ws.Rows("8").RowHeight = 28.8
$ws.Range("B8").Value = 7

$ws.Range("A9").Style = "Normal"
$ws.Range("A9").Value = "provide Sabisu authentication in the application"
$ws.Range("B9").Value = 5

# Row 10: "add QR code scanning " (unchanged text/style), days now 5.
$ws.Range("B10").Value = 5

# Match the final active selection cell recorded in the workbook.
$ws.Range("A8").Select()
